$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: date value increments by 1 day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Update prices in column D for rows 29-36
$ws.Range("D29").Value = 1619
$ws.Range("D30").Value = 1896
$ws.Range("D31").Value = 1223
$ws.Range("D32").Value = 1365
$ws.Range("D33").Value = 1326
$ws.Range("D34").Value = 1577
$ws.Range("D35").Value = 1521
$ws.Range("D36").Value = 1779
